{"js": "// Office.js (Word JavaScript API) implementation of the diff:\n//  1. Insert a new, un-styled paragraph at the very top of the body containing\n//     the literal text \"<1T>\" and carrying the \"_GoBack\" bookmark (id 0).\n//  2. Remove the \"_GoBack\" bookmark from its old spot (middle of the final\n//     paragraph) and re-join that paragraph's text into a single run.\n//  3. On the four \"Heading1\" definition paragraphs (s.update,\n//     s.intersection_update, s.difference_update,\n//     s.symmetric_difference_update) drop the <w:proofErr w:type=\"gramStart\"/>\n//     / <w:proofErr w:type=\"gramEnd\"/> markers and merge the split function-name\n//     runs back into one run (spellStart/spellEnd stay).\n//\n// Strategy: use Range.insertOoxml(\"Replace\") (flat-OPC WordprocessingML) to\n// splice in exact OOXML for each paragraph that needs new shape. This avoids\n// relying on higher level formatting APIs that can't express proofErr /\n// bookmark placement, while leaving every untouched paragraph byte-for-byte\n// alone.\n\nconst NS = 'xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"';\n\nfunction flatOpc(bodyInnerXml) {\n  return (\n    '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document ' + NS + '><w:body>' + bodyInnerXml + '</w:body></w:document></pkg:xmlData>' +\n    '</pkg:part></pkg:package>'\n  );\n}\n\n// --- 1. New first paragraph with the marker text + relocated bookmark ----\nconst body = context.document.body;\nconst markerPara = body.insertParagraph(\"\", \"Start\");\nawait context.sync();\n\nconst markerXml = flatOpc(\n  '<w:p>' +\n    '<w:r><w:t>&lt;1T&gt;</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>'\n);\nmarkerPara.getRange(\"Whole\").insertOoxml(markerXml, \"Replace\");\nawait context.sync();\n\n// --- 2. The four Heading1 \"definition\" paragraphs: drop gramStart/gramEnd,\n//        merge the split name runs into one run -------------------------\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nfunction headingXml(name, rest) {\n  return flatOpc(\n    '<w:p>' +\n      '<w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:r><w:t>' + name + '</w:t></w:r>' +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t>' + rest + '</w:t></w:r>' +\n    '</w:p>'\n  );\n}\n\nconst headingFixes = [\n  { needle: \"s.update(*others)\", name: \"s.update\", rest: \"(*others) [alias: s |= other | \\u2026]\" },\n  { needle: \"s.intersection_update(*others)\", name: \"s.intersection_update\", rest: \"(*others) [alias: s &amp;= other &amp; \\u2026]\" },\n  { needle: \"s.difference_update(*others)\", name: \"s.difference_update\", rest: \"(*others) [alias: s -= other | \\u2026]\" },\n  { needle: \"s.symmetric_difference_update(other)\", name: \"s.symmetric_difference_update\", rest: \"(other) [alias: s ^= other]\" },\n];\n\nfor (const fix of headingFixes) {\n  const items = paras.items;\n  let target = null;\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(fix.needle) === 0) {\n      target = items[i];\n      break;\n    }\n  }\n  if (!target) {\n    throw new Error(\"Could not locate heading paragraph for: \" + fix.needle);\n  }\n  target.getRange(\"Whole\").insertOoxml(headingXml(fix.name, fix.rest), \"Replace\");\n  await context.sync();\n  paras.load(\"items/text\");\n  await context.sync();\n}\n\n// --- 3. Last paragraph: remove the stale bookmark, rejoin the text -------\nconst finalXml = flatOpc(\n  '<w:p>' +\n    '<w:r><w:tab/><w:t>Update the set s, keeping only elements found in either set, but not in both</w:t></w:r>' +\n  '</w:p>'\n);\n\nconst items = paras.items;\nlet lastTarget = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"but not in both\") > -1) {\n    lastTarget = items[i];\n  }\n}\nif (!lastTarget) {\n  // fall back: the final paragraph in the body\n  lastTarget = items[items.length - 1];\n}\nlastTarget.getRange(\"Whole\").insertOoxml(finalXml, \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop implementation of the diff:\n#  1. Insert a new, un-styled paragraph at the very top of the body containing\n#     the literal text \"<1T>\" and carrying the \"_GoBack\" bookmark (id 0).\n#  2. Remove the \"_GoBack\" bookmark from its old spot (middle of the final\n#     paragraph) and re-join that paragraph's text into a single run.\n#  3. On the four \"Heading1\" definition paragraphs (s.update,\n#     s.intersection_update, s.difference_update,\n#     s.symmetric_difference_update) drop the <w:proofErr w:type=\"gramStart\"/>\n#     / <w:proofErr w:type=\"gramEnd\"/> markers and merge the split function-name\n#     runs back into one run (spellStart/spellEnd stay).\n#\n# Strategy: use Range.InsertXML (flat-OPC WordprocessingML), the native COM\n# counterpart of Office.js's Range.insertOoxml, to splice exact OOXML into\n# each paragraph that needs a new shape, leaving every untouched paragraph\n# byte-for-byte alone.\n\n$d = $word.ActiveDocument\n\nfunction FlatOpc($bodyInnerXml) {\n    return '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n        $bodyInnerXml +\n        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# --- 1. New first paragraph with the marker text + relocated bookmark ----\n$markerXml = FlatOpc('<w:p><w:r><w:t>&lt;1T&gt;</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>')\n$head = $d.Range(0, 0)\n$head.Collapse(1)\n$head.InsertXML($markerXml)\n\n# --- 2. The four Heading1 \"definition\" paragraphs: drop gramStart/gramEnd,\n#        merge the split name runs into one run -------------------------\nfunction HeadingXml($name, $rest) {\n    return FlatOpc('<w:p><w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr><w:proofErr w:type=\"spellStart\"/><w:r><w:t>' + $name + '</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>' + $rest + '</w:t></w:r></w:p>')\n}\n\n$ellipsis = [char]0x2026\n\n$headingFixes = @(\n    @{ Needle = \"s.update(\"; Name = \"s.update\"; Rest = (\"(*others) [alias: s |= other | \" + $ellipsis + \"]\") },\n    @{ Needle = \"s.intersection_update(\"; Name = \"s.intersection_update\"; Rest = (\"(*others) [alias: s &amp;= other &amp; \" + $ellipsis + \"]\") },\n    @{ Needle = \"s.difference_update(\"; Name = \"s.difference_update\"; Rest = (\"(*others) [alias: s -= other | \" + $ellipsis + \"]\") },\n    @{ Needle = \"s.symmetric_difference_update(\"; Name = \"s.symmetric_difference_update\"; Rest = \"(other) [alias: s ^= other]\" }\n)\n\nforeach ($fix in $headingFixes) {\n    $found = $null\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Text.StartsWith($fix.Needle)) {\n            $found = $p\n            break\n        }\n    }\n    if ($found -eq $null) {\n        throw (\"Could not locate heading paragraph for: \" + $fix.Needle)\n    }\n    $xml = HeadingXml $fix.Name $fix.Rest\n    $found.Range.InsertXML($xml)\n}\n\n# --- 3. Last paragraph: remove the stale bookmark, rejoin the text -------\n# Trim the trailing paragraph mark off the range so InsertXML replaces the\n# paragraph's content in place instead of leaving a spare empty paragraph\n# behind (the body's very last w:p can't simply be \"deleted\").\n$finalXml = FlatOpc('<w:p><w:r><w:tab/><w:t>Update the set s, keeping only elements found in either set, but not in both</w:t></w:r></w:p>')\n\n$lastIndex = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($lastIndex)\n$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)\n$lastRange.InsertXML($finalXml)\n"}
